# Auto-generated Excel COM-interop script to apply the diff changes
# to Sheets/Aegis_Profits.xlsx (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 129
$ws.Range("H129").Value = 3559.4167
$ws.Range("J129").Value = 1018.8947
$ws.Range("L129").Value = 3056.6841
$ws.Range("N129").Value = -13056.6841

# Row 132
$ws.Range("H132").Value = 9267749
$ws.Range("I132").Value = 10008769
$ws.Range("K132").Value = 30026307
$ws.Range("M132").Value = -30023777

# Row 138
$ws.Range("H138").Value = 2238.753
$ws.Range("I138").Value = 1762
$ws.Range("J138").Value = 2343.2466
$ws.Range("K138").Value = 5286
$ws.Range("L138").Value = 7029.739799999999
$ws.Range("M138").Value = -146
$ws.Range("N138").Value = -17309.7398

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 96914
$ws.Range("I45").Value = 168716.83
$ws.Range("J45").Value = 10750.6
$ws.Range("K45").Value = 168716.83
$ws.Range("L45").Value = 10750.6
$ws.Range("M45").Value = -168339.83
$ws.Range("N45").Value = -11504.6

# Row 61
$ws.Range("H61").Value = 2435.7273
$ws.Range("I61").Value = 1757
$ws.Range("J61").Value = 3623.5
$ws.Range("K61").Value = 1757
$ws.Range("L61").Value = 3623.5
$ws.Range("M61").Value = -1545
$ws.Range("N61").Value = -4047.5

# Row 110
$ws.Range("H110").Value = 125263130
$ws.Range("I110").Value = 125263130
$ws.Range("K110").Value = 125263130
$ws.Range("M110").Value = -125261085

# Row 122
$ws.Range("H122").Value = 1992.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1992.5
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 5977.5
$ws.Range("N122").Value = -10877.5

# Row 132
$ws.Range("H132").Value = 19773.354
$ws.Range("I132").Value = 29342.523
$ws.Range("J132").Value = 4315.4614
$ws.Range("K132").Value = 88027.569
$ws.Range("L132").Value = 12946.3842
$ws.Range("M132").Value = -85497.569
$ws.Range("N132").Value = -18006.3842

# Row 136
$ws.Range("H136").Value = 2435.7273
$ws.Range("I136").Value = 1757
$ws.Range("J136").Value = 3623.5
$ws.Range("K136").Value = 5271
$ws.Range("L136").Value = 10870.5
$ws.Range("M136").Value = -2721
$ws.Range("N136").Value = -15970.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 735.1667
$ws.Range("I94").Value = 702.0625
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 702.0625
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -251.0625
$ws.Range("N94").Value = -1902

# Row 105
$ws.Range("H105").Value = 92785.45
$ws.Range("I105").Value = 68583.266
$ws.Range("J105").Value = 144647.28
$ws.Range("K105").Value = 68583.266
$ws.Range("L105").Value = 144647.28
$ws.Range("M105").Value = -66836.266
$ws.Range("N105").Value = -148141.28

# Row 107
$ws.Range("H107").Value = 90909976
$ws.Range("I107").Value = 111111970
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 111111970
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -111110050
$ws.Range("N107").Value = -4840

# Row 134
$ws.Range("H134").Value = 26127.38
$ws.Range("I134").Value = 34547.133
$ws.Range("J134").Value = 5078
$ws.Range("K134").Value = 103641.399
$ws.Range("L134").Value = 15234
$ws.Range("M134").Value = -101106.399
$ws.Range("N134").Value = -20304

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1666.3334
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 28
$ws.Range("H28").Value = 28088.6
$ws.Range("J28").Value = 28088.6
$ws.Range("L28").Value = 28088.6
$ws.Range("N28").Value = -28578.6

# Row 99
$ws.Range("H99").Value = 13117.889
$ws.Range("I99").Value = 2443.3333
$ws.Range("J99").Value = 18455.166
$ws.Range("K99").Value = 2443.3333
$ws.Range("L99").Value = 18455.166
$ws.Range("M99").Value = -945.3332999999998
$ws.Range("N99").Value = -21451.166

# Row 107
$ws.Range("H107").Value = 1010.4
$ws.Range("J107").Value = 975
$ws.Range("L107").Value = 975
$ws.Range("N107").Value = -4815

# Row 113
$ws.Range("H113").Value = 1666.3334
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# Row 126
$ws.Range("H126").Value = 13117.889
$ws.Range("I126").Value = 2443.3333
$ws.Range("J126").Value = 18455.166
$ws.Range("K126").Value = 7329.999899999999
$ws.Range("L126").Value = 55365.49800000001
$ws.Range("M126").Value = -4859.999899999999
$ws.Range("N126").Value = -60305.49800000001

$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 4089
$ws.Range("I56").Value = 4089
$ws.Range("K56").Value = 4089
$ws.Range("M56").Value = -3559

# Row 137
$ws.Range("H137").Value = 38241.07
$ws.Range("I137").Value = 68392.664
$ws.Range("J137").Value = 3450.7693
$ws.Range("K137").Value = 205177.992
$ws.Range("L137").Value = 10352.3079
$ws.Range("M137").Value = -200077.992
$ws.Range("N137").Value = -20552.3079

$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 9000
$ws.Range("J52").Value = 9000
$ws.Range("L52").Value = 9000
$ws.Range("N52").Value = -9518

# Row 80
$ws.Range("H80").Value = 55614430
$ws.Range("I80").Value = 111226960
$ws.Range("J80").Value = 1901.4445
$ws.Range("K80").Value = 111226960
$ws.Range("L80").Value = 1901.4445
$ws.Range("M80").Value = -111225962
$ws.Range("N80").Value = -3897.4445

# Row 83
$ws.Range("H83").Value = 55614430
$ws.Range("I83").Value = 111226960
$ws.Range("J83").Value = 1901.4445
$ws.Range("K83").Value = 556134800
$ws.Range("L83").Value = 9507.2225
$ws.Range("M83").Value = -556129808
$ws.Range("N83").Value = -19491.2225

# Row 107
$ws.Range("H107").Value = 631802.5600000001
$ws.Range("I107").Value = 328.2
$ws.Range("J107").Value = 1684259.9
$ws.Range("K107").Value = 328.2
$ws.Range("L107").Value = 1684259.9
$ws.Range("M107").Value = 1591.8
$ws.Range("N107").Value = -1688099.9

# Row 113
$ws.Range("H113").Value = 2034.1875
$ws.Range("I113").Value = 2217.5715
$ws.Range("J113").Value = 1891.5555
$ws.Range("K113").Value = 2217.5715
$ws.Range("L113").Value = 1891.5555
$ws.Range("M113").Value = -47.57150000000001
$ws.Range("N113").Value = -6231.5555

# Row 122
$ws.Range("H122").Value = 2650
$ws.Range("I122").Value = 2742.8572
$ws.Range("K122").Value = 8228.571599999999
$ws.Range("M122").Value = -5778.571599999999

# Row 132
$ws.Range("H132").Value = 2566.7908
$ws.Range("I132").Value = 2023.3823
$ws.Range("K132").Value = 6070.1469
$ws.Range("M132").Value = -3540.1469

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1685.5555
$ws.Range("I100").Value = 1400
$ws.Range("J100").Value = 2134.2856
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 2134.2856
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -3216.2856

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2587.2942
$ws.Range("I122").Value = 1865.2667
$ws.Range("K122").Value = 5595.800099999999
$ws.Range("M122").Value = -3145.800099999999

# Row 132
$ws.Range("H132").Value = 2822.465
$ws.Range("I132").Value = 2895.3823
$ws.Range("J132").Value = 2547
$ws.Range("K132").Value = 8686.1469
$ws.Range("L132").Value = 7641
$ws.Range("M132").Value = -6156.1469
$ws.Range("N132").Value = -12701
